# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# For each changed row we update Price (col D) and/or Volume(1h) (col E); a few
# rows (32/34 and 42/43) also swap which coin occupies that rank, so Coin (B)
# and Link (C) are rewritten too.
#
# NumberFormat is forced to "@" (Text) before writing numeric-looking Price
# strings (e.g. "1.00", "8.60") so COM's type inference doesn't collapse them
# into plain numbers and silently drop significant trailing zeros - matching
# the source data, which stores every Price/Volume cell as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "69.798.21"
$ws.Cells.Item(2, 5).Value = "  +0.48%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.707.86"
$ws.Cells.Item(3, 5).Value = "  +0.49%  "

$ws.Cells.Item(4, 5).Value = "  -0.03%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "678.32"
$ws.Cells.Item(5, 5).Value = "  -1.14%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "162.65"
$ws.Cells.Item(6, 5).Value = "  +1.58%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.00"
$ws.Cells.Item(7, 5).Value = "  -0.04%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.498"
$ws.Cells.Item(8, 5).Value = "  +0.78%  "

$ws.Cells.Item(9, 5).Value = "  +1.73%  "

$ws.Cells.Item(10, 5).Value = "  +0.92%  "

$ws.Cells.Item(11, 5).Value = "  +2.17%  "

$ws.Cells.Item(12, 5).Value = "  +1.21%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "32.97"
$ws.Cells.Item(13, 5).Value = "  +1.29%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "3.718.18"
$ws.Cells.Item(14, 5).Value = "  +0.35%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "69.783.00"
$ws.Cells.Item(15, 5).Value = "  +0.53%  "

$ws.Cells.Item(16, 5).Value = "  +1.77%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "16.15"
$ws.Cells.Item(17, 5).Value = "  +1.86%  "

$ws.Cells.Item(18, 5).Value = "  +1.55%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "473.97"
$ws.Cells.Item(19, 5).Value = "  +0.55%  "

$ws.Cells.Item(20, 5).Value = "  -1.57%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.655"
$ws.Cells.Item(21, 5).Value = "  +0.53%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "80.54"
$ws.Cells.Item(22, 5).Value = "  +0.97%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "3.854.50"
$ws.Cells.Item(23, 5).Value = "  +0.47%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.0000129"
$ws.Cells.Item(24, 5).Value = "  +3.78%  "

$ws.Cells.Item(25, 5).Value = "  +0.02%  "

$ws.Cells.Item(26, 5).Value = "  +0.05%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "9.17"
$ws.Cells.Item(27, 5).Value = "  -0.96%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.72"
$ws.Cells.Item(28, 5).Value = "  -0.22%  "

$ws.Cells.Item(29, 5).Value = "  +0.60%  "

$ws.Cells.Item(30, 5).Value = "  +1.25%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "6.64"
$ws.Cells.Item(31, 5).Value = "  +0.87%  "

$ws.Cells.Item(32, 2).Value = "EthereumClassic"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "27.02"
$ws.Cells.Item(32, 5).Value = "  +0.29%  "

$ws.Cells.Item(33, 5).Value = "  +0.23%  "

$ws.Cells.Item(34, 2).Value = "Kaspa"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.166"
$ws.Cells.Item(34, 5).Value = "  +4.02%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.696.65"
$ws.Cells.Item(35, 5).Value = "  +0.92%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "8.60"
$ws.Cells.Item(36, 5).Value = "  +4.79%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "6.21"
$ws.Cells.Item(37, 5).Value = "  +0.87%  "

$ws.Cells.Item(39, 5).Value = "  +0.41%  "

$ws.Cells.Item(40, 5).Value = "  -0.02%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0909"
$ws.Cells.Item(41, 5).Value = "  +1.06%  "

$ws.Cells.Item(42, 2).Value = "Monero"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "167.58"
$ws.Cells.Item(42, 5).Value = "  +1.09%  "

$ws.Cells.Item(43, 2).Value = "Mantle"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.947"
$ws.Cells.Item(43, 5).Value = "  +0.44%  "

$ws.Cells.Item(44, 5).Value = "  -1.16%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.81"
$ws.Cells.Item(45, 5).Value = "  +2.59%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.000283"
$ws.Cells.Item(46, 5).Value = "  -0.35%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "28.25"
$ws.Cells.Item(47, 5).Value = "  +1.13%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.12"
$ws.Cells.Item(48, 5).Value = "  -1.10%  "

$ws.Cells.Item(49, 5).Value = "  +0.09%  "

$ws.Cells.Item(50, 5).Value = "  +2.05%  "

$ws.Cells.Item(51, 5).Value = "  +2.46%  "
